# Auto-generated edit script for Jerseys.xlsx (Jersey Vitals + Jersey Colors sheets)
$wb = $excel.ActiveWorkbook

# ---------- Sheet 1: "Jersey Vitals" ----------
$ws1 = $wb.Worksheets.Item("Jersey Vitals")

# Remove the two trailing rows (old rows 27 "Uniform File" and 28 "UNIQUEID")
# so the sheet shrinks from 28 rows to 26 rows, matching the new layout.
$ws1.Range("A27:H28").EntireRow.Delete() | Out-Null

# Rewrite rows 2-26 with the updated offsets / normalized (UPPER_SNAKE_CASE) names.
$ws1.Cells.Item(2, 1).Value = "COLOR_LUMINANCE_LEVEL"
$ws1.Cells.Item(2, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(2, 3).Value = "COLORLUMINANCELEVEL"
$ws1.Cells.Item(2, 4).Value = ""
$ws1.Cells.Item(2, 5).Value = "0x143 (type=slider, length=32, startBit=4)"
$ws1.Cells.Item(2, 6).Value = "0x143 (type=slider, length=32, startBit=4)"
$ws1.Cells.Item(2, 7).Value = ""
$ws1.Cells.Item(2, 8).Value = ""

$ws1.Cells.Item(3, 1).Value = "EDITION_CRC32B"
$ws1.Cells.Item(3, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(3, 3).Value = "EDITIONCRC32B"
$ws1.Cells.Item(3, 4).Value = "0x54 (type=combo, length=32, startBit=0)"
$ws1.Cells.Item(3, 5).Value = "0x54 (type=combo, length=32, startBit=0)"
$ws1.Cells.Item(3, 6).Value = "0x54 (type=combo, length=32, startBit=0)"
$ws1.Cells.Item(3, 7).Value = ""
$ws1.Cells.Item(3, 8).Value = ""

$ws1.Cells.Item(4, 1).Value = "EDITION_NAME"
$ws1.Cells.Item(4, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(4, 3).Value = "EDITIONNAME"
$ws1.Cells.Item(4, 4).Value = "0x5C (type=string, length=512, startBit=0)"
$ws1.Cells.Item(4, 5).Value = "0x5C (type=string, length=512, startBit=0)"
$ws1.Cells.Item(4, 6).Value = "0x5C (type=string, length=512, startBit=0)"
$ws1.Cells.Item(4, 7).Value = ""
$ws1.Cells.Item(4, 8).Value = ""

$ws1.Cells.Item(5, 1).Value = "FILENAME"
$ws1.Cells.Item(5, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(5, 3).Value = "FILENAME"
$ws1.Cells.Item(5, 4).Value = "0x8 (type=string, length=512, startBit=0)"
$ws1.Cells.Item(5, 5).Value = "0x8 (type=string, length=512, startBit=0)"
$ws1.Cells.Item(5, 6).Value = "0x8 (type=string, length=512, startBit=0)"
$ws1.Cells.Item(5, 7).Value = ""
$ws1.Cells.Item(5, 8).Value = ""

$ws1.Cells.Item(6, 1).Value = "HEADBAND_LOGO_TYPE"
$ws1.Cells.Item(6, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(6, 3).Value = "HEADBANDLOGOTYPE"
$ws1.Cells.Item(6, 4).Value = "0x14E (type=combo, length=2, startBit=3)"
$ws1.Cells.Item(6, 5).Value = "0x14E (type=combo, length=2, startBit=3)"
$ws1.Cells.Item(6, 6).Value = "0x14E (type=combo, length=2, startBit=3)"
$ws1.Cells.Item(6, 7).Value = ""
$ws1.Cells.Item(6, 8).Value = ""

$ws1.Cells.Item(7, 1).Value = "IS_ALTERNATE"
$ws1.Cells.Item(7, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(7, 3).Value = "ISALTERNATE"
$ws1.Cells.Item(7, 4).Value = "0x14A (type=combo, length=1, startBit=5)"
$ws1.Cells.Item(7, 5).Value = "0x14A (type=combo, length=1, startBit=5)"
$ws1.Cells.Item(7, 6).Value = "0x14A (type=combo, length=1, startBit=5)"
$ws1.Cells.Item(7, 7).Value = ""
$ws1.Cells.Item(7, 8).Value = ""

$ws1.Cells.Item(8, 1).Value = "IS_HOME"
$ws1.Cells.Item(8, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(8, 3).Value = "ISHOME"
$ws1.Cells.Item(8, 4).Value = ""
$ws1.Cells.Item(8, 5).Value = "0x14A (type=combo, length=1, startBit=0)"
$ws1.Cells.Item(8, 6).Value = "0x14A (type=combo, length=1, startBit=1)"
$ws1.Cells.Item(8, 7).Value = ""
$ws1.Cells.Item(8, 8).Value = ""

$ws1.Cells.Item(9, 1).Value = "IS_TEAM_CREATE_UNIFORM"
$ws1.Cells.Item(9, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(9, 3).Value = "ISTEAMCREATEUNIFORM"
$ws1.Cells.Item(9, 4).Value = "0x14A (type=combo, length=1, startBit=1)"
$ws1.Cells.Item(9, 5).Value = "0x14A (type=combo, length=1, startBit=1)"
$ws1.Cells.Item(9, 6).Value = "0x14A (type=combo, length=1, startBit=2)"
$ws1.Cells.Item(9, 7).Value = ""
$ws1.Cells.Item(9, 8).Value = ""

$ws1.Cells.Item(10, 1).Value = "JACKET_WARMUP_CRC32B"
$ws1.Cells.Item(10, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(10, 3).Value = "JACKETWARMUPCRC32B"
$ws1.Cells.Item(10, 4).Value = "0xDC (type=number, length=32, startBit=0)"
$ws1.Cells.Item(10, 5).Value = "0xDC (type=number, length=32, startBit=0)"
$ws1.Cells.Item(10, 6).Value = "0xDC (type=number, length=32, startBit=0)"
$ws1.Cells.Item(10, 7).Value = ""
$ws1.Cells.Item(10, 8).Value = ""

$ws1.Cells.Item(11, 1).Value = "LOGO_BRAND"
$ws1.Cells.Item(11, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(11, 3).Value = "LOGOBRAND"
$ws1.Cells.Item(11, 4).Value = "0x143 (type=combo, length=2, startBit=0)"
$ws1.Cells.Item(11, 5).Value = "0x143 (type=combo, length=2, startBit=0)"
$ws1.Cells.Item(11, 6).Value = "0x143 (type=combo, length=2, startBit=0)"
$ws1.Cells.Item(11, 7).Value = ""
$ws1.Cells.Item(11, 8).Value = ""

$ws1.Cells.Item(12, 1).Value = "LOGO_TYPE"
$ws1.Cells.Item(12, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(12, 3).Value = "LOGOTYPE"
$ws1.Cells.Item(12, 4).Value = "0x142 (type=combo, length=3, startBit=4)"
$ws1.Cells.Item(12, 5).Value = "0x142 (type=combo, length=3, startBit=4)"
$ws1.Cells.Item(12, 6).Value = "0x142 (type=combo, length=3, startBit=4)"
$ws1.Cells.Item(12, 7).Value = ""
$ws1.Cells.Item(12, 8).Value = ""

$ws1.Cells.Item(13, 1).Value = "MYTEAM_INCLUDE"
$ws1.Cells.Item(13, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(13, 3).Value = "MYTEAMINCLUDE"
$ws1.Cells.Item(13, 4).Value = "0x14A (type=combo, length=1, startBit=3)"
$ws1.Cells.Item(13, 5).Value = "0x14A (type=combo, length=1, startBit=3)"
$ws1.Cells.Item(13, 6).Value = "0x14A (type=combo, length=1, startBit=3)"
$ws1.Cells.Item(13, 7).Value = ""
$ws1.Cells.Item(13, 8).Value = ""

$ws1.Cells.Item(14, 1).Value = "NUMBER_ON_SHORTS"
$ws1.Cells.Item(14, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(14, 3).Value = "NUMBERONSHORTS"
$ws1.Cells.Item(14, 4).Value = "0x14A (type=combo, length=1, startBit=4)"
$ws1.Cells.Item(14, 5).Value = "0x14A (type=combo, length=1, startBit=4)"
$ws1.Cells.Item(14, 6).Value = "0x14A (type=combo, length=1, startBit=4)"
$ws1.Cells.Item(14, 7).Value = ""
$ws1.Cells.Item(14, 8).Value = ""

$ws1.Cells.Item(15, 1).Value = "OVERRIDE_NBA_PATCH"
$ws1.Cells.Item(15, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(15, 3).Value = "OVERRIDENBAPATCH"
$ws1.Cells.Item(15, 4).Value = "0x14A (type=combo, length=1, startBit=2)"
$ws1.Cells.Item(15, 5).Value = "0x14A (type=combo, length=1, startBit=2)"
$ws1.Cells.Item(15, 6).Value = "0x14A (type=combo, length=1, startBit=2)"
$ws1.Cells.Item(15, 7).Value = ""
$ws1.Cells.Item(15, 8).Value = ""

$ws1.Cells.Item(16, 1).Value = "SOCKS_AWAY"
$ws1.Cells.Item(16, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(16, 3).Value = "SOCKSAWAY"
$ws1.Cells.Item(16, 4).Value = "0xD4 (type=combo, length=32, startBit=0)"
$ws1.Cells.Item(16, 5).Value = "0xD4 (type=combo, length=32, startBit=0)"
$ws1.Cells.Item(16, 6).Value = "0xD4 (type=combo, length=32, startBit=0)"
$ws1.Cells.Item(16, 7).Value = ""
$ws1.Cells.Item(16, 8).Value = ""

$ws1.Cells.Item(17, 1).Value = "SOCKS_COLOR_N#1"
$ws1.Cells.Item(17, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(17, 3).Value = "SOCKSCOLORN1"
$ws1.Cells.Item(17, 4).Value = ""
$ws1.Cells.Item(17, 5).Value = "0xE0 (type=combo, length=32, startBit=0)"
$ws1.Cells.Item(17, 6).Value = "0xE0 (type=combo, length=32, startBit=0)"
$ws1.Cells.Item(17, 7).Value = ""
$ws1.Cells.Item(17, 8).Value = ""

$ws1.Cells.Item(18, 1).Value = "SOCKS_COLOR_N#2"
$ws1.Cells.Item(18, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(18, 3).Value = "SOCKSCOLORN2"
$ws1.Cells.Item(18, 4).Value = ""
$ws1.Cells.Item(18, 5).Value = "0xE4 (type=combo, length=32, startBit=0)"
$ws1.Cells.Item(18, 6).Value = "0xE4 (type=combo, length=32, startBit=0)"
$ws1.Cells.Item(18, 7).Value = ""
$ws1.Cells.Item(18, 8).Value = ""

$ws1.Cells.Item(19, 1).Value = "SOCKS_COLOR_N#3"
$ws1.Cells.Item(19, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(19, 3).Value = "SOCKSCOLORN3"
$ws1.Cells.Item(19, 4).Value = ""
$ws1.Cells.Item(19, 5).Value = "0xE8 (type=combo, length=32, startBit=0)"
$ws1.Cells.Item(19, 6).Value = "0xE8 (type=combo, length=32, startBit=0)"
$ws1.Cells.Item(19, 7).Value = ""
$ws1.Cells.Item(19, 8).Value = ""

$ws1.Cells.Item(20, 1).Value = "SOCKS_COLOR_N#4"
$ws1.Cells.Item(20, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(20, 3).Value = "SOCKSCOLORN4"
$ws1.Cells.Item(20, 4).Value = ""
$ws1.Cells.Item(20, 5).Value = "0xEC (type=combo, length=32, startBit=0)"
$ws1.Cells.Item(20, 6).Value = "0xEC (type=combo, length=32, startBit=0)"
$ws1.Cells.Item(20, 7).Value = ""
$ws1.Cells.Item(20, 8).Value = ""

$ws1.Cells.Item(21, 1).Value = "SOCKS_COLOR_N#5"
$ws1.Cells.Item(21, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(21, 3).Value = "SOCKSCOLORN5"
$ws1.Cells.Item(21, 4).Value = ""
$ws1.Cells.Item(21, 5).Value = "0xF0 (type=combo, length=32, startBit=0)"
$ws1.Cells.Item(21, 6).Value = "0xF0 (type=combo, length=32, startBit=0)"
$ws1.Cells.Item(21, 7).Value = ""
$ws1.Cells.Item(21, 8).Value = ""

$ws1.Cells.Item(22, 1).Value = "SOCKS_HOME"
$ws1.Cells.Item(22, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(22, 3).Value = "SOCKSHOME"
$ws1.Cells.Item(22, 4).Value = "0xD0 (type=combo, length=32, startBit=0)"
$ws1.Cells.Item(22, 5).Value = "0xD0 (type=combo, length=32, startBit=0)"
$ws1.Cells.Item(22, 6).Value = "0xD0 (type=combo, length=32, startBit=0)"
$ws1.Cells.Item(22, 7).Value = ""
$ws1.Cells.Item(22, 8).Value = ""

$ws1.Cells.Item(23, 1).Value = "SPONSOR_PATCH"
$ws1.Cells.Item(23, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(23, 3).Value = "SPONSORPATCH"
$ws1.Cells.Item(23, 4).Value = "0x48 (type=string, length=64, startBit=0)"
$ws1.Cells.Item(23, 5).Value = "0x48 (type=string, length=64, startBit=0)"
$ws1.Cells.Item(23, 6).Value = "0x48 (type=string, length=64, startBit=0)"
$ws1.Cells.Item(23, 7).Value = ""
$ws1.Cells.Item(23, 8).Value = ""

$ws1.Cells.Item(24, 1).Value = "TEAM"
$ws1.Cells.Item(24, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(24, 3).Value = "TEAM"
$ws1.Cells.Item(24, 4).Value = "0x52 (type=combo, length=16, startBit=0)"
$ws1.Cells.Item(24, 5).Value = "0x52 (type=combo, length=16, startBit=0)"
$ws1.Cells.Item(24, 6).Value = "0x52 (type=combo, length=16, startBit=0)"
$ws1.Cells.Item(24, 7).Value = ""
$ws1.Cells.Item(24, 8).Value = ""

$ws1.Cells.Item(25, 1).Value = "TYPE"
$ws1.Cells.Item(25, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(25, 3).Value = "TYPE"
$ws1.Cells.Item(25, 4).Value = "0x140 (type=combo, length=5, startBit=0)"
$ws1.Cells.Item(25, 5).Value = "0x140 (type=combo, length=5, startBit=0)"
$ws1.Cells.Item(25, 6).Value = "0x140 (type=combo, length=5, startBit=0)"
$ws1.Cells.Item(25, 7).Value = ""
$ws1.Cells.Item(25, 8).Value = ""

$ws1.Cells.Item(26, 1).Value = "UNIQUEID"
$ws1.Cells.Item(26, 2).Value = "Jersey Vitals"
$ws1.Cells.Item(26, 3).Value = "UNIQUEID"
$ws1.Cells.Item(26, 4).Value = "0x50 (type=number, length=16, startBit=0)"
$ws1.Cells.Item(26, 5).Value = "0x50 (type=number, length=16, startBit=0)"
$ws1.Cells.Item(26, 6).Value = "0x50 (type=number, length=16, startBit=0)"
$ws1.Cells.Item(26, 7).Value = ""
$ws1.Cells.Item(26, 8).Value = ""

# ---------- Sheet 2: "Jersey Colors" ----------
$ws2 = $wb.Worksheets.Item("Jersey Colors")

# Rewrite rows 2-28 (sheet grows from 24 to 28 rows: offsets reshuffled and four new
# accessory/sock color variants are introduced: PRIMARY_COLOR#ACCESSORY, PRIMARY_COLOR#SOCK,
# QUATERNARY_COLOR#ACCESSORY, QUATERNARY_COLOR#SOCK, QUINTARY_COLOR#ACCESSORY, QUINTARY_COLOR#SOCK,
# SECONDARY_COLOR#ACCESSORY, SECONDARY_COLOR#SOCK, TERTIARY_COLOR#ACCESSORY, TERTIARY_COLOR#SOCK.
$ws2.Cells.Item(2, 1).Value = "ARM_ACCESSORY_COLOR_AWAY"
$ws2.Cells.Item(2, 2).Value = "Jersey Colors"
$ws2.Cells.Item(2, 3).Value = "ARMACCESSORYCOLORAWAY"
$ws2.Cells.Item(2, 4).Value = "0x14D (type=combo, length=3, startBit=1)"
$ws2.Cells.Item(2, 5).Value = "0x14D (type=combo, length=3, startBit=1)"
$ws2.Cells.Item(2, 6).Value = "0x14D (type=combo, length=3, startBit=1)"
$ws2.Cells.Item(2, 7).Value = ""
$ws2.Cells.Item(2, 8).Value = ""

$ws2.Cells.Item(3, 1).Value = "ARM_ACCESSORY_COLOR_HOME"
$ws2.Cells.Item(3, 2).Value = "Jersey Colors"
$ws2.Cells.Item(3, 3).Value = "ARMACCESSORYCOLORHOME"
$ws2.Cells.Item(3, 4).Value = "0x14C (type=combo, length=3, startBit=6)"
$ws2.Cells.Item(3, 5).Value = "0x14C (type=combo, length=3, startBit=6)"
$ws2.Cells.Item(3, 6).Value = "0x14C (type=combo, length=3, startBit=6)"
$ws2.Cells.Item(3, 7).Value = ""
$ws2.Cells.Item(3, 8).Value = ""

$ws2.Cells.Item(4, 1).Value = "HEADBAND_COLOR_AWAY"
$ws2.Cells.Item(4, 2).Value = "Jersey Colors"
$ws2.Cells.Item(4, 3).Value = "HEADBANDCOLORAWAY"
$ws2.Cells.Item(4, 4).Value = "0x14D (type=combo, length=3, startBit=7)"
$ws2.Cells.Item(4, 5).Value = "0x14D (type=combo, length=3, startBit=7)"
$ws2.Cells.Item(4, 6).Value = "0x14D (type=combo, length=3, startBit=7)"
$ws2.Cells.Item(4, 7).Value = ""
$ws2.Cells.Item(4, 8).Value = ""

$ws2.Cells.Item(5, 1).Value = "HEADBAND_COLOR_HOME"
$ws2.Cells.Item(5, 2).Value = "Jersey Colors"
$ws2.Cells.Item(5, 3).Value = "HEADBANDCOLORHOME"
$ws2.Cells.Item(5, 4).Value = "0x14D (type=combo, length=3, startBit=4)"
$ws2.Cells.Item(5, 5).Value = "0x14D (type=combo, length=3, startBit=4)"
$ws2.Cells.Item(5, 6).Value = "0x14D (type=combo, length=3, startBit=4)"
$ws2.Cells.Item(5, 7).Value = ""
$ws2.Cells.Item(5, 8).Value = ""

$ws2.Cells.Item(6, 1).Value = "LEG_ACCESSORY_COLOR_AWAY"
$ws2.Cells.Item(6, 2).Value = "Jersey Colors"
$ws2.Cells.Item(6, 3).Value = "LEGACCESSORYCOLORAWAY"
$ws2.Cells.Item(6, 4).Value = "0x14C (type=combo, length=3, startBit=3)"
$ws2.Cells.Item(6, 5).Value = "0x14C (type=combo, length=3, startBit=3)"
$ws2.Cells.Item(6, 6).Value = "0x14C (type=combo, length=3, startBit=3)"
$ws2.Cells.Item(6, 7).Value = ""
$ws2.Cells.Item(6, 8).Value = ""

$ws2.Cells.Item(7, 1).Value = "LEG_ACCESSORY_COLOR_HOME"
$ws2.Cells.Item(7, 2).Value = "Jersey Colors"
$ws2.Cells.Item(7, 3).Value = "LEGACCESSORYCOLORHOME"
$ws2.Cells.Item(7, 4).Value = "0x14C (type=combo, length=3, startBit=0)"
$ws2.Cells.Item(7, 5).Value = "0x14C (type=combo, length=3, startBit=0)"
$ws2.Cells.Item(7, 6).Value = "0x14C (type=combo, length=3, startBit=0)"
$ws2.Cells.Item(7, 7).Value = ""
$ws2.Cells.Item(7, 8).Value = ""

$ws2.Cells.Item(8, 1).Value = "PRIMARY_COLOR"
$ws2.Cells.Item(8, 2).Value = "Jersey Colors"
$ws2.Cells.Item(8, 3).Value = "PRIMARYCOLOR"
$ws2.Cells.Item(8, 4).Value = "0xF4 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(8, 5).Value = "0xF4 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(8, 6).Value = "0xF4 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(8, 7).Value = ""
$ws2.Cells.Item(8, 8).Value = ""

$ws2.Cells.Item(9, 1).Value = "PRIMARY_COLOR#ACCESSORY"
$ws2.Cells.Item(9, 2).Value = "Jersey Colors"
$ws2.Cells.Item(9, 3).Value = "PRIMARYCOLORACCESSORY"
$ws2.Cells.Item(9, 4).Value = "0x10C (type=color, length=24, startBit=0)"
$ws2.Cells.Item(9, 5).Value = "0x10C (type=color, length=24, startBit=0)"
$ws2.Cells.Item(9, 6).Value = "0x10C (type=color, length=24, startBit=0)"
$ws2.Cells.Item(9, 7).Value = ""
$ws2.Cells.Item(9, 8).Value = ""

$ws2.Cells.Item(10, 1).Value = "PRIMARY_COLOR#SOCK"
$ws2.Cells.Item(10, 2).Value = "Jersey Colors"
$ws2.Cells.Item(10, 3).Value = "PRIMARYCOLORSOCK"
$ws2.Cells.Item(10, 4).Value = "0x12C (type=color, length=24, startBit=0)"
$ws2.Cells.Item(10, 5).Value = "0x12C (type=color, length=24, startBit=0)"
$ws2.Cells.Item(10, 6).Value = "0x12C (type=color, length=24, startBit=0)"
$ws2.Cells.Item(10, 7).Value = ""
$ws2.Cells.Item(10, 8).Value = ""

$ws2.Cells.Item(11, 1).Value = "QUATERNARY_COLOR"
$ws2.Cells.Item(11, 2).Value = "Jersey Colors"
$ws2.Cells.Item(11, 3).Value = "QUATERNARYCOLOR"
$ws2.Cells.Item(11, 4).Value = "0x100 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(11, 5).Value = "0x100 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(11, 6).Value = "0x100 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(11, 7).Value = ""
$ws2.Cells.Item(11, 8).Value = ""

$ws2.Cells.Item(12, 1).Value = "QUATERNARY_COLOR#ACCESSORY"
$ws2.Cells.Item(12, 2).Value = "Jersey Colors"
$ws2.Cells.Item(12, 3).Value = "QUATERNARYCOLORACCESSORY"
$ws2.Cells.Item(12, 4).Value = "0x118 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(12, 5).Value = "0x118 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(12, 6).Value = "0x118 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(12, 7).Value = ""
$ws2.Cells.Item(12, 8).Value = ""

$ws2.Cells.Item(13, 1).Value = "QUATERNARY_COLOR#SOCK"
$ws2.Cells.Item(13, 2).Value = "Jersey Colors"
$ws2.Cells.Item(13, 3).Value = "QUATERNARYCOLORSOCK"
$ws2.Cells.Item(13, 4).Value = "0x138 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(13, 5).Value = "0x138 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(13, 6).Value = "0x138 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(13, 7).Value = ""
$ws2.Cells.Item(13, 8).Value = ""

$ws2.Cells.Item(14, 1).Value = "QUINTARY_COLOR"
$ws2.Cells.Item(14, 2).Value = "Jersey Colors"
$ws2.Cells.Item(14, 3).Value = "QUINTARYCOLOR"
$ws2.Cells.Item(14, 4).Value = "0x104 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(14, 5).Value = "0x104 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(14, 6).Value = "0x104 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(14, 7).Value = ""
$ws2.Cells.Item(14, 8).Value = ""

$ws2.Cells.Item(15, 1).Value = "QUINTARY_COLOR#ACCESSORY"
$ws2.Cells.Item(15, 2).Value = "Jersey Colors"
$ws2.Cells.Item(15, 3).Value = "QUINTARYCOLORACCESSORY"
$ws2.Cells.Item(15, 4).Value = "0x11C (type=color, length=24, startBit=0)"
$ws2.Cells.Item(15, 5).Value = "0x11C (type=color, length=24, startBit=0)"
$ws2.Cells.Item(15, 6).Value = "0x11C (type=color, length=24, startBit=0)"
$ws2.Cells.Item(15, 7).Value = ""
$ws2.Cells.Item(15, 8).Value = ""

$ws2.Cells.Item(16, 1).Value = "QUINTARY_COLOR#SOCK"
$ws2.Cells.Item(16, 2).Value = "Jersey Colors"
$ws2.Cells.Item(16, 3).Value = "QUINTARYCOLORSOCK"
$ws2.Cells.Item(16, 4).Value = "0x13C (type=color, length=24, startBit=0)"
$ws2.Cells.Item(16, 5).Value = "0x13C (type=color, length=24, startBit=0)"
$ws2.Cells.Item(16, 6).Value = "0x13C (type=color, length=24, startBit=0)"
$ws2.Cells.Item(16, 7).Value = ""
$ws2.Cells.Item(16, 8).Value = ""

$ws2.Cells.Item(17, 1).Value = "SECONDARY_COLOR#ACCESSORY"
$ws2.Cells.Item(17, 2).Value = "Jersey Colors"
$ws2.Cells.Item(17, 3).Value = "SECONDARYCOLORACCESSORY"
$ws2.Cells.Item(17, 4).Value = "0x110 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(17, 5).Value = "0x110 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(17, 6).Value = "0x110 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(17, 7).Value = ""
$ws2.Cells.Item(17, 8).Value = ""

$ws2.Cells.Item(18, 1).Value = "SECONDARY_COLOR#SOCK"
$ws2.Cells.Item(18, 2).Value = "Jersey Colors"
$ws2.Cells.Item(18, 3).Value = "SECONDARYCOLORSOCK"
$ws2.Cells.Item(18, 4).Value = "0x130 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(18, 5).Value = "0x130 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(18, 6).Value = "0x130 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(18, 7).Value = ""
$ws2.Cells.Item(18, 8).Value = ""

$ws2.Cells.Item(19, 1).Value = "SECONDARY_COLOR_SHOES_COLOR"
$ws2.Cells.Item(19, 2).Value = "Jersey Colors"
$ws2.Cells.Item(19, 3).Value = "SECONDARYCOLORSHOESCOLOR"
$ws2.Cells.Item(19, 4).Value = "0xF8 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(19, 5).Value = "0xF8 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(19, 6).Value = "0xF8 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(19, 7).Value = ""
$ws2.Cells.Item(19, 8).Value = ""

$ws2.Cells.Item(20, 1).Value = "SEXTARY_COLOR"
$ws2.Cells.Item(20, 2).Value = "Jersey Colors"
$ws2.Cells.Item(20, 3).Value = "SEXTARYCOLOR"
$ws2.Cells.Item(20, 4).Value = "0x108 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(20, 5).Value = "0x108 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(20, 6).Value = "0x108 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(20, 7).Value = ""
$ws2.Cells.Item(20, 8).Value = ""

$ws2.Cells.Item(21, 1).Value = "SHOE_PRIMARY_COLOR_AWAY"
$ws2.Cells.Item(21, 2).Value = "Jersey Colors"
$ws2.Cells.Item(21, 3).Value = "SHOEPRIMARYCOLORAWAY"
$ws2.Cells.Item(21, 4).Value = "0x128 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(21, 5).Value = "0x128 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(21, 6).Value = "0x128 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(21, 7).Value = ""
$ws2.Cells.Item(21, 8).Value = ""

$ws2.Cells.Item(22, 1).Value = "SOCK_COLOR_AWAY"
$ws2.Cells.Item(22, 2).Value = "Jersey Colors"
$ws2.Cells.Item(22, 3).Value = "SOCKCOLORAWAY"
$ws2.Cells.Item(22, 4).Value = "0x14F (type=combo, length=3, startBit=0)"
$ws2.Cells.Item(22, 5).Value = "0x14F (type=combo, length=3, startBit=0)"
$ws2.Cells.Item(22, 6).Value = "0x14F (type=combo, length=3, startBit=0)"
$ws2.Cells.Item(22, 7).Value = ""
$ws2.Cells.Item(22, 8).Value = ""

$ws2.Cells.Item(23, 1).Value = "SOCK_COLOR_HOME"
$ws2.Cells.Item(23, 2).Value = "Jersey Colors"
$ws2.Cells.Item(23, 3).Value = "SOCKCOLORHOME"
$ws2.Cells.Item(23, 4).Value = "0x14E (type=combo, length=3, startBit=5)"
$ws2.Cells.Item(23, 5).Value = "0x14E (type=combo, length=3, startBit=5)"
$ws2.Cells.Item(23, 6).Value = "0x14E (type=combo, length=3, startBit=5)"
$ws2.Cells.Item(23, 7).Value = ""
$ws2.Cells.Item(23, 8).Value = ""

$ws2.Cells.Item(24, 1).Value = "TERTIARY_COLOR"
$ws2.Cells.Item(24, 2).Value = "Jersey Colors"
$ws2.Cells.Item(24, 3).Value = "TERTIARYCOLOR"
$ws2.Cells.Item(24, 4).Value = "0xFC (type=color, length=24, startBit=0)"
$ws2.Cells.Item(24, 5).Value = "0xFC (type=color, length=24, startBit=0)"
$ws2.Cells.Item(24, 6).Value = "0xFC (type=color, length=24, startBit=0)"
$ws2.Cells.Item(24, 7).Value = ""
$ws2.Cells.Item(24, 8).Value = ""

$ws2.Cells.Item(25, 1).Value = "TERTIARY_COLOR#ACCESSORY"
$ws2.Cells.Item(25, 2).Value = "Jersey Colors"
$ws2.Cells.Item(25, 3).Value = "TERTIARYCOLORACCESSORY"
$ws2.Cells.Item(25, 4).Value = "0x114 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(25, 5).Value = "0x114 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(25, 6).Value = "0x114 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(25, 7).Value = ""
$ws2.Cells.Item(25, 8).Value = ""

$ws2.Cells.Item(26, 1).Value = "TERTIARY_COLOR#SOCK"
$ws2.Cells.Item(26, 2).Value = "Jersey Colors"
$ws2.Cells.Item(26, 3).Value = "TERTIARYCOLORSOCK"
$ws2.Cells.Item(26, 4).Value = "0x134 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(26, 5).Value = "0x134 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(26, 6).Value = "0x134 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(26, 7).Value = ""
$ws2.Cells.Item(26, 8).Value = ""

$ws2.Cells.Item(27, 1).Value = "UNDER_SHIRT_COLOR_AWAY"
$ws2.Cells.Item(27, 2).Value = "Jersey Colors"
$ws2.Cells.Item(27, 3).Value = "UNDERSHIRTCOLORAWAY"
$ws2.Cells.Item(27, 4).Value = "0x124 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(27, 5).Value = "0x124 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(27, 6).Value = "0x124 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(27, 7).Value = ""
$ws2.Cells.Item(27, 8).Value = ""

$ws2.Cells.Item(28, 1).Value = "UNDER_SHIRT_COLOR_HOME"
$ws2.Cells.Item(28, 2).Value = "Jersey Colors"
$ws2.Cells.Item(28, 3).Value = "UNDERSHIRTCOLORHOME"
$ws2.Cells.Item(28, 4).Value = "0x120 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(28, 5).Value = "0x120 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(28, 6).Value = "0x120 (type=color, length=24, startBit=0)"
$ws2.Cells.Item(28, 7).Value = ""
$ws2.Cells.Item(28, 8).Value = ""

